$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (E1/F1/G1). Shared strings must land in this order:
#   48 = actgr_PCcarsNum (F1)
#   49 = actgr_CCcarsNum (E1)
#   50 = actgr_CCcarsRent (G1)
$ws.Range("F1").Value = "actgr_PCcarsNum"
$ws.Range("E1").Value = "actgr_CCcarsNum"
$ws.Range("G1").Value = "actgr_CCcarsRent"

# New numeric data for columns E (CCcarsNum), F (PCcarsNum), G (CCcarsRent)
$data = @(
    @{ Row = 2;  E = 2; F = 1;  G = 171000 },
    @{ Row = 3;  E = 4; F = 3;  G = 173700 },
    @{ Row = 4;  E = 1; F = 1;  G = 175200 },
    @{ Row = 5;  E = 4; F = 6;  G = 174480 },
    @{ Row = 6;  E = 9; F = 10; G = 156401.32999999999 },
    @{ Row = 7;  E = 3; F = 5;  G = 175600 },
    @{ Row = 8;  E = 4; F = 4;  G = 159000 },
    @{ Row = 9;  E = 3; F = 2;  G = 170400 },
    @{ Row = 10; E = 4; F = 7;  G = 166950 },
    @{ Row = 11; E = 2; F = 1;  G = 164100 },
    @{ Row = 12; E = 2; F = 2;  G = 171600 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
}

# New number format (numFmtId 2, "0.00") on the rent column
$ws.Range("G2:G12").NumberFormat = "0.00"

# Column widths: A narrows, E/F/G get new bestFit-style widths.
# (Values chosen land on the closest pixel the COM width model can represent
# to the target 5.375 / 16.5 / 16.25 / 16.5 character widths.)
$ws.Columns.Item(1).ColumnWidth = 4.714285714285714
$ws.Columns.Item(5).ColumnWidth = 15.714285714285714
$ws.Columns.Item(6).ColumnWidth = 15.571428571428571
$ws.Columns.Item(7).ColumnWidth = 15.714285714285714

# Selection moves to the new G1 header cell
$ws.Range("G1").Select() | Out-Null
